# This script reproduces the target diff:
#   1) "M . Raja"      -> "R ."  (proofErr gramStart/gramEnd) + " Akshitha"
#   2) "BATCH : 36"     -> proofErr gramStart + "BATCH" + " :" (proofErr gramEnd) + " 36"
#   3) "2303A5" + "2277" (two runs) -> single run "2303A52464"
#
# Plain Find/Replace on this runtime auto-merges adjacent runs that share
# identical formatting, which makes it impossible to reproduce the
# w:proofErr-delimited run splits shown in the diff. Instead we target the
# exact Range covering the text to change and call Range.InsertXML with a
# <w:p>-wrapped OOXML fragment: the runtime then substitutes just that
# range's run content in place (run splits, proofErr markers and all)
# without disturbing the surrounding document structure.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Change 1 (NAME line): "M . Raja" -> "R ." / proofErr / " Akshitha"
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("M . Raja")
if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + 8)
    $xml = '<w:p ' + $wNs + '>' +
             '<w:proofErr w:type="gramStart"/>' +
             '<w:r><w:rPr><w:b/><w:spacing w:val="-2"/></w:rPr><w:t>R .</w:t></w:r>' +
             '<w:proofErr w:type="gramEnd"/>' +
             '<w:r><w:rPr><w:b/><w:spacing w:val="-2"/></w:rPr><w:t xml:space="preserve"> Akshitha</w:t></w:r>' +
           '</w:p>'
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Change 2 (BATCH line): "BATCH : 36" -> proofErr / "BATCH" / " :" /
#                         proofErr / " 36"
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("BATCH : 36")
if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + 10)
    $xml = '<w:p ' + $wNs + '>' +
             '<w:proofErr w:type="gramStart"/>' +
             '<w:r><w:rPr><w:b/><w:spacing w:val="-2"/></w:rPr><w:t>BATCH</w:t></w:r>' +
             '<w:r><w:rPr><w:b/><w:spacing w:val="-2"/></w:rPr><w:t xml:space="preserve"> :</w:t></w:r>' +
             '<w:proofErr w:type="gramEnd"/>' +
             '<w:r><w:rPr><w:b/><w:spacing w:val="-2"/></w:rPr><w:t xml:space="preserve"> 36</w:t></w:r>' +
           '</w:p>'
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Change 3 (HALLTICKET.NO line): "2303A5" + "2277" (two runs) -> single
#                                  run "2303A52464"
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("2303A52277")
if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + 10)
    $xml = '<w:p ' + $wNs + '>' +
             '<w:r><w:rPr><w:b/><w:spacing w:val="-4"/></w:rPr><w:t>2303A52464</w:t></w:r>' +
           '</w:p>'
    $r.InsertXML($xml)
}

Write-Host "Edits applied."
